$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1231.3158
$ws.Range("J98").Value = 8498
$ws.Range("L98").Value = 8498
$ws.Range("N98").Value = -11494
$ws.Range("H115").Value = 2305.4
$ws.Range("I115").Value = 2845.6667
$ws.Range("K115").Value = 8537.000100000001
$ws.Range("M115").Value = -6970.000100000001
$ws.Range("H122").Value = 1231.3158
$ws.Range("J122").Value = 8498
$ws.Range("L122").Value = 25494
$ws.Range("N122").Value = -30394
$ws.Range("H135").Value = 72194.86
$ws.Range("I135").Value = 658.1667
$ws.Range("J135").Value = 125847.375
$ws.Range("K135").Value = 5923.5003
$ws.Range("L135").Value = 1132626.375
$ws.Range("M135").Value = -3388.5003
$ws.Range("N135").Value = -1137696.375
$ws.Range("H138").Value = 4577.654
$ws.Range("I138").Value = 8149.5
$ws.Range("J138").Value = 3928.2273
$ws.Range("K138").Value = 24448.5
$ws.Range("L138").Value = 11784.6819
$ws.Range("M138").Value = -19308.5
$ws.Range("N138").Value = -22064.6819
$ws.Range("H141").Value = 2063.8958
$ws.Range("I141").Value = 1755.7693
$ws.Range("J141").Value = 3399.111
$ws.Range("K141").Value = 5267.3079
$ws.Range("L141").Value = 10197.333
$ws.Range("M141").Value = -87.30789999999979
$ws.Range("N141").Value = -20557.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28579.219
$ws.Range("I32").Value = 36012.5
$ws.Range("J32").Value = 8306.637000000001
$ws.Range("K32").Value = 36012.5
$ws.Range("L32").Value = 8306.637000000001
$ws.Range("M32").Value = -35725.5
$ws.Range("N32").Value = -8880.637000000001
$ws.Range("H45").Value = 2513.2856
$ws.Range("I45").Value = 2044.65
$ws.Range("K45").Value = 2044.65
$ws.Range("M45").Value = -1667.65
$ws.Range("H88").Value = 2583.9524
$ws.Range("I88").Value = 1277.6
$ws.Range("J88").Value = 2992.1875
$ws.Range("K88").Value = 1277.6
$ws.Range("L88").Value = 2992.1875
$ws.Range("M88").Value = -871.5999999999999
$ws.Range("N88").Value = -3804.1875
$ws.Range("H91").Value = 2583.9524
$ws.Range("I91").Value = 1277.6
$ws.Range("J91").Value = 2992.1875
$ws.Range("K91").Value = 1277.6
$ws.Range("L91").Value = 2992.1875
$ws.Range("M91").Value = 126.4000000000001
$ws.Range("N91").Value = -5800.1875
$ws.Range("H102").Value = 1732.2572
$ws.Range("I102").Value = 1524.3667
$ws.Range("K102").Value = 1524.3667
$ws.Range("M102").Value = 97.63329999999996
$ws.Range("H132").Value = 29155.684
$ws.Range("I132").Value = 35987.766
$ws.Range("J132").Value = 3535.375
$ws.Range("K132").Value = 107963.298
$ws.Range("L132").Value = 10606.125
$ws.Range("M132").Value = -105433.298
$ws.Range("N132").Value = -15666.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 22501.5
$ws.Range("I7").Value = 22501.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 22501.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -22388.5
$ws.Range("N7").ClearContents()
$ws.Range("H86").Value = 3132
$ws.Range("I86").Value = 2007.9
$ws.Range("K86").Value = 2007.9
$ws.Range("M86").Value = -884.9000000000001
$ws.Range("H89").Value = 3132
$ws.Range("I89").Value = 2007.9
$ws.Range("K89").Value = 10039.5
$ws.Range("M89").Value = -4423.5
$ws.Range("H99").Value = 50563.41
$ws.Range("I99").Value = 66040.44
$ws.Range("J99").Value = 9291.333000000001
$ws.Range("K99").Value = 66040.44
$ws.Range("L99").Value = 9291.333000000001
$ws.Range("M99").Value = -64542.44
$ws.Range("N99").Value = -12287.333
$ws.Range("H105").Value = 3335.923
$ws.Range("I105").Value = 3446.8
$ws.Range("K105").Value = 3446.8
$ws.Range("M105").Value = -1699.8
$ws.Range("H107").Value = 2221
$ws.Range("I107").Value = 1566.2727
$ws.Range("K107").Value = 1566.2727
$ws.Range("M107").Value = 353.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1163.9
$ws.Range("I6").Value = 846.6667
$ws.Range("J6").Value = 1299.8572
$ws.Range("K6").Value = 846.6667
$ws.Range("L6").Value = 1299.8572
$ws.Range("M6").Value = -733.6667
$ws.Range("N6").Value = -1525.8572
$ws.Range("H12").Value = 6869.6665
$ws.Range("I12").Value = 4804.75
$ws.Range("J12").Value = 10999.5
$ws.Range("K12").Value = 4804.75
$ws.Range("L12").Value = 10999.5
$ws.Range("M12").Value = -4634.75
$ws.Range("N12").Value = -11339.5
$ws.Range("H21").Value = 4500
$ws.Range("J21").Value = 4500
$ws.Range("L21").Value = 4500
$ws.Range("N21").Value = -4970
$ws.Range("H31").Value = 3384.875
$ws.Range("J31").Value = 3488.4546
$ws.Range("L31").Value = 3488.4546
$ws.Range("N31").Value = -4078.4546
$ws.Range("H34").Value = 3384.875
$ws.Range("J34").Value = 3488.4546
$ws.Range("L34").Value = 3488.4546
$ws.Range("N34").Value = -3892.4546
$ws.Range("H105").Value = 979.7273
$ws.Range("I105").Value = 979.7273
$ws.Range("K105").Value = 979.7273
$ws.Range("M105").Value = 767.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2074.5454
$ws.Range("J68").Value = 1936.3334
$ws.Range("L68").Value = 5809.0002
$ws.Range("N68").Value = -7431.0002
$ws.Range("H71").Value = 2074.5454
$ws.Range("J71").Value = 1936.3334
$ws.Range("L71").Value = 17427.0006
$ws.Range("N71").Value = -25539.0006
$ws.Range("H113").Value = 1273
$ws.Range("I113").Value = 740.5714
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2221.7142
$ws.Range("L113").Value = 15000
$ws.Range("M113").Value = -51.71420000000035
$ws.Range("N113").Value = -19340
$ws.Range("H140").Value = 2438.7917
$ws.Range("I140").Value = 2587.1904
$ws.Range("J140").Value = 1400
$ws.Range("K140").Value = 7761.5712
$ws.Range("L140").Value = 4200
$ws.Range("M140").Value = -2581.5712
$ws.Range("N140").Value = -14560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3890.6191
$ws.Range("I70").Value = 3570.9412
$ws.Range("J70").Value = 5249.25
$ws.Range("K70").Value = 3570.9412
$ws.Range("L70").Value = 5249.25
$ws.Range("M70").Value = -3300.9412
$ws.Range("N70").Value = -5789.25
$ws.Range("H73").Value = 3890.6191
$ws.Range("I73").Value = 3570.9412
$ws.Range("J73").Value = 5249.25
$ws.Range("K73").Value = 3570.9412
$ws.Range("L73").Value = 5249.25
$ws.Range("M73").Value = -2634.9412
$ws.Range("N73").Value = -7121.25
$ws.Range("H80").Value = 5887
$ws.Range("I80").Value = 2396
$ws.Range("K80").Value = 2396
$ws.Range("M80").Value = -1398
$ws.Range("H83").Value = 5887
$ws.Range("I83").Value = 2396
$ws.Range("K83").Value = 11980
$ws.Range("M83").Value = -6988
$ws.Range("H132").Value = 42807.12
$ws.Range("I132").Value = 44424.168
$ws.Range("J132").Value = 3998
$ws.Range("K132").Value = 133272.504
$ws.Range("L132").Value = 11994
$ws.Range("M132").Value = -130742.504
$ws.Range("N132").Value = -17054

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 19.5
$ws.Range("I4").Value = 19.5
$ws.Range("K4").Value = 19.5
$ws.Range("M4").Value = 93.5
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H17").Value = 2656
$ws.Range("I17").Value = 2656
$ws.Range("K17").Value = 2656
$ws.Range("M17").Value = -2484
$ws.Range("H136").Value = 3696.0232
$ws.Range("I136").Value = 2696.7917
$ws.Range("J136").Value = 4958.2104
$ws.Range("K136").Value = 8090.375100000001
$ws.Range("L136").Value = 14874.6312
$ws.Range("M136").Value = -5540.375100000001
$ws.Range("N136").Value = -19974.6312
